$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "helix jump"
$ws.Range("B7").Value = "com.singleton.helix"

$ws.Range("A8").Select()
